$wb = $excel.ActiveWorkbook

# --- Sheet1: add a new "policy" header cell in A1, move selection to D1 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1").Value = "policy"
$ws1.Range("D1").Select()

# --- Sheet3: change the I2 parameter value (10000 -> 10), update selection,
#     and make Sheet3 the active sheet/tab ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()
$ws3.Range("I2").Value = 10
$ws3.Range("I2").Select()
